$wb = $excel.ActiveWorkbook

# --- Rename sheets (tab order matches workbook.xml sheet order) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511686885946965"
$wb.Worksheets.Item(2).Name = "NB_TO-16511686913865025"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168691387472"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511686914334712"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651168691511471"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168688562694.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686885766938.csv"
$ws1.Range("B4").Value = "go_stims-16511686885776973.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686885937119.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511686913704693.csv"
$ws2.Range("B3").Value = "ZB-match_7-16511686886627052.csv"
$ws2.Range("B4").Value = "OB-16511686899744716.csv"
$ws2.Range("B5").Value = "OB-16511686898584852.csv"
$ws2.Range("B6").Value = "TB-16511686912255077.csv"
$ws2.Range("B7").Value = "OB-1651168689409471.csv"
$ws2.Range("B8").Value = "ZB-match_9-16511686886246958.csv"
$ws2.Range("B9").Value = "TB-16511686900805037.csv"
$ws2.Range("B10").Value = "ZB-match_6-1651168688777696.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511686914014692.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686913894727.csv"
$ws4.Range("B4").Value = "MM_stims-16511686914164698.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686914014692.csv"
$ws4.Range("B6").Value = "MM_stims-16511686914324746.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168691417471.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511686914954693.csv"
$ws5.Range("B3").Value = "SAT_stims-16511686914634798.csv"
$ws5.Range("B4").Value = "SAT_stims-1651168691437472.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651168691479473.csv"
